$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror row 12 (date + description text, with its formatting) into a new row 13
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial()
